$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 10 change from serial date 45208 to 45212
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
